# Remove the duplicate "description" column (column O) from the bacteria
# metadata sheet. Column F already holds "description"; column O was a
# duplicate header that should be deleted entirely, shifting all the
# following columns (P.. / case info / assembly info / SRA fields /
# test fields) one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Activate()
$ws.Columns("O").Delete()

# Restore a sensible view/selection on the sheet after the edit.
$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
